# Update the "Förändrad" (changed) date column (C) for all data rows
# from 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
